$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-12-12 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-13 Wednesday", 2) | Out-Null

# Update each answer cell in the table by position, to avoid ambiguity
# from values that coincidentally match other cells after editing.
$t = $d.Tables.Item(1)

# Row 1, Col 1: "19÷4=4, 3" -> "71÷7=10, 1"
$t.Cell(1, 1).Range.Text = "71÷7=10, 1"

# Row 1, Col 2: "44÷9=4, 8" -> "44÷2=22, 0"
$t.Cell(1, 2).Range.Text = "44÷2=22, 0"

# Row 1, Col 3: "18÷9=2, 0" -> "25÷6=4, 1"
$t.Cell(1, 3).Range.Text = "25÷6=4, 1"

# Row 1, Col 4: "99÷3=33, 0" -> "47÷7=6, 5"
$t.Cell(1, 4).Range.Text = "47÷7=6, 5"

# Row 1, Col 5: "65÷2=32, 1" -> "58÷7=8, 2"
$t.Cell(1, 5).Range.Text = "58÷7=8, 2"

# Row 5, Col 1: "73÷4=18, 1" -> "32÷4=8, 0"
$t.Cell(5, 1).Range.Text = "32÷4=8, 0"

# Row 5, Col 2: "92÷3=30, 2" -> "41÷5=8, 1"
$t.Cell(5, 2).Range.Text = "41÷5=8, 1"

# Row 5, Col 3: "86÷2=43, 0" -> "47÷7=6, 5"
$t.Cell(5, 3).Range.Text = "47÷7=6, 5"

# Row 5, Col 4: "37÷3=12, 1" -> "43÷6=7, 1"
$t.Cell(5, 4).Range.Text = "43÷6=7, 1"

# Row 5, Col 5: "83÷5=16, 3" -> "59÷8=7, 3"
$t.Cell(5, 5).Range.Text = "59÷8=7, 3"

# Row 9, Col 1: "16÷8=2, 0" -> "14÷5=2, 4"
$t.Cell(9, 1).Range.Text = "14÷5=2, 4"

# Row 9, Col 2: "74÷9=8, 2" -> "40÷4=10, 0"
$t.Cell(9, 2).Range.Text = "40÷4=10, 0"

# Row 9, Col 3: "74÷8=9, 2" -> "39÷7=5, 4"
$t.Cell(9, 3).Range.Text = "39÷7=5, 4"

# Row 9, Col 4: "60÷4=15, 0" -> "12÷5=2, 2"
$t.Cell(9, 4).Range.Text = "12÷5=2, 2"

# Row 9, Col 5: "16÷2=8, 0" -> "19÷3=6, 1"
$t.Cell(9, 5).Range.Text = "19÷3=6, 1"

# Row 13, Col 1: "32÷8=4, 0" -> "70÷7=10, 0"
$t.Cell(13, 1).Range.Text = "70÷7=10, 0"

# Row 13, Col 2: "84÷8=10, 4" -> "24÷3=8, 0"
$t.Cell(13, 2).Range.Text = "24÷3=8, 0"

# Row 13, Col 3: "35÷6=5, 5" -> "28÷3=9, 1"
$t.Cell(13, 3).Range.Text = "28÷3=9, 1"

# Row 13, Col 4: "29÷2=14, 1" -> "81÷2=40, 1"
$t.Cell(13, 4).Range.Text = "81÷2=40, 1"

# Row 13, Col 5: "28÷3=9, 1" -> "67÷6=11, 1"
$t.Cell(13, 5).Range.Text = "67÷6=11, 1"

# Row 17, Col 1: "69÷3=23, 0" -> "42÷8=5, 2"
$t.Cell(17, 1).Range.Text = "42÷8=5, 2"

# Row 17, Col 2: "22÷9=2, 4" -> "59÷7=8, 3"
$t.Cell(17, 2).Range.Text = "59÷7=8, 3"

# Row 17, Col 3: "70÷7=10, 0" -> "94÷9=10, 4"
$t.Cell(17, 3).Range.Text = "94÷9=10, 4"

# Row 17, Col 4: "82÷6=13, 4" -> "36÷9=4, 0"
$t.Cell(17, 4).Range.Text = "36÷9=4, 0"

# Row 17, Col 5: "65÷9=7, 2" -> "88÷8=11, 0"
$t.Cell(17, 5).Range.Text = "88÷8=11, 0"

